$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "27-12-2024"
$ws.Range("B2").Value = "20:00"

# Force C2 to hold the text string "2" (not a number) without leaving any
# extra cell-style (quotePrefix / number-format) behind: compute it via a
# formula that yields a text result, then paste only the value back.
$ws.Range("Z1").Formula = '="2"'
$ws.Range("Z1").Copy()
$ws.Range("C2").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

$ws.Range("D2").Value = "preuba"
$ws.Range("E2").Value = "ai"
$ws.Range("F2").Value = 25000
